# Added correction Cls and Cds
# New columns L:N ("Dynamic Pressure (Corrected)", "CL (corrected)", "CD (corrected)")
# on the "Tail On 65 MPH" sheet, computed from eps_tot_t.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tail On 65 MPH")

# Header row
$ws.Range("L1").Value = "Dynamic Pressure (Corrected)"
$ws.Range("M1").Value = "CL (corrected)"
$ws.Range("N1").Value = "CD (corrected)"

# Data rows 2:85 - write formulas cell-by-cell (not as one range) so each
# row keeps its own full formula text instead of being compressed into a
# shared formula.
for ($r = 2; $r -le 85; $r++) {
    $ws.Range("L$r").Formula = "=0.5*(65*12)^2*0.002377*(1+eps_tot_t)^2"
    $ws.Range("M$r").Formula = "=B$r/(L$r*wing_area)"
    $ws.Range("N$r").Formula = "=D$r/(L$r*wing_area)"
}

# Number format on the corrected dynamic-pressure column
$ws.Range("L2:L85").NumberFormat = "0.0"

# Column widths for the newly-added columns
$ws.Columns.Item(12).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(13).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(14).EntireColumn.AutoFit() | Out-Null

# Update view/selection state: "Tail On 65 MPH" becomes the active sheet
# (it was "Airplane Characteristics"), with the new O1 cell selected.
$ws5 = $wb.Worksheets.Item("Airplane Characteristics")
$ws5.Range("G34").Select() | Out-Null

$ws.Select() | Out-Null
$ws.Range("O1").Select() | Out-Null
